$d = $word.ActiveDocument

# Replace the run-text content of a paragraph (identified by its current
# text) with freshly authored OOXML, while leaving the paragraph's own
# mark / pPr / any other sibling runs (e.g. a leading empty <w:r/>) intact.
# We shrink the paragraph's range by one character (the trailing paragraph
# mark) before calling InsertXML so the paragraph mark itself is never
# disturbed - this avoids spurious extra paragraphs when the edited
# paragraph happens to be the last one in the body.
function Set-ParagraphRunXml($matchText, $runXml) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$matchText*") {
            $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
            $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
                   '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body><w:p>' + $runXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $rng.InsertXML($pkg)
            return
        }
    }
}

# Title heading (H1) - shorten it.
Set-ParagraphRunXml 'Play Game of Cards for Free - Review of Online Slot Game' `
    '<w:r><w:t>Play Game of Cards for Free</w:t></w:r>'

# "What we like" bullet list.
Set-ParagraphRunXml 'Exciting theme based on "Alice in Wonderland"' `
    '<w:r><w:t>Exciting theme based on Alice in Wonderland</w:t></w:r>'

Set-ParagraphRunXml '3 Bonus features for more chances to win big' `
    '<w:r><w:t>Multiple Bonus features for increased winning potential</w:t></w:r>'

Set-ParagraphRunXml 'Wide range of betting options' `
    '<w:r><w:t>Wide range of betting levels to suit different players</w:t></w:r>'

Set-ParagraphRunXml 'Progressive Jackpot available' `
    '<w:r><w:t>High jackpot multiplier of 2,000x line bet</w:t></w:r>'

# "What we don't like" bullet list.
Set-ParagraphRunXml 'Only 30 fixed paylines' `
    '<w:r><w:t>Medium volatility may not appeal to players seeking high-risk/high-reward gameplay</w:t></w:r>'

Set-ParagraphRunXml 'Medium volatility may not suit those seeking high-risk options' `
    '<w:r><w:t>Limited number of fixed paylines (30) compared to some other slots</w:t></w:r>'

# Closing bold title line.
Set-ParagraphRunXml 'Play Game of Cards for Free - Review of Online Slot Game' `
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Game of Cards for Free</w:t></w:r>'

# Closing italic summary line.
Set-ParagraphRunXml 'Join the battle between the Queen of Hearts and the Queen of Spades in Game of Cards. Play for free and read our review of this exciting online slot game.' `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Game of Cards and play for free to experience the exciting Alice in Wonderland-themed slot game.</w:t></w:r>'
